$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Dummy"
$ws.Range("A3").Value = 0

$ws.Range("F3").Select()
